$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: was an empty cell with "no-border" style (s=2); now gets a value
# and picks up the plain text style (s=1), matching the pattern used by the
# surrounding rows (e.g. A44). Filled first so the new shared string it
# introduces ("D184T22YLD28") lands at index 6.
$ws.Range("A44").Copy()
$ws.Range("A43").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A43").Value = "D184T22YLD28"

# Row 39: same treatment, filled second so its new shared string
# ("D184T21AAAAA") lands at index 7.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A39").Value = "D184T21AAAAA"

$excel.CutCopyMode = $false

# Move the active selection from A37 to A39, as recorded in the saved view.
[void]$ws.Range("A39").Select()
